$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the method signature used in the second example block (row 25)
$ws.Range("C25").Value = "SmartRules MyDatatype[] myRule(MyDatatype inputParam, String x)"

# The array-returning rule result now yields "= null" (quote-prefixed text) instead of numbers
$ws.Range("D27").Value = "'= null"
$ws.Range("D28").Value = "'= null"
$ws.Range("D29").Value = "'= null"

# Update the view state: scrolled so row 7 is the top-left row, and D29 selected
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D29").Select()
